# voto_score_matrix.xlsx: "Gemeinsame Liste" is replaced by "Freie Wähler" (same column/row slot,
# recomputed scores) and a brand-new party "Tierschutzpartei" is inserted, growing the symmetric
# score matrix from 15x15 (A1:O15) to 16x16 (A1:P16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert one column at M and one row at 13 for "Tierschutzpartei", which shifts the
# trailing UFFBASSE / Volt / WGD columns and rows one slot to the right / down.
$ws.Range("M1").EntireColumn.Insert() | Out-Null
$ws.Range("A13").EntireRow.Insert() | Out-Null

# Write the full, refreshed 16x16 score matrix (header row/column + symmetric pairwise scores).
# Row 1
$ws.Range("A1").Value = "p1"
$ws.Range("B1").Value = "AfD"
$ws.Range("C1").Value = "BSW"
$ws.Range("D1").Value = "CDU"
$ws.Range("E1").Value = "DaGe"
$ws.Range("F1").Value = "Die Grünen"
$ws.Range("G1").Value = "Die Partei"
$ws.Range("H1").Value = "FDP"
$ws.Range("I1").Value = "Freie Wähler"
$ws.Range("J1").Value = "Linke"
$ws.Range("K1").Value = "SPD"
$ws.Range("L1").Value = "SfD"
$ws.Range("M1").Value = "Tierschutzpartei"
$ws.Range("N1").Value = "UFFBASSE"
$ws.Range("O1").Value = "Volt"
$ws.Range("P1").Value = "WGD"
# Row 2
$ws.Range("A2").Value = "AfD"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.4451219512195122
$ws.Range("D2").Value = 0.5914634146341463
$ws.Range("E2").Value = 0.3170731707317073
$ws.Range("F2").Value = 0.2347560975609756
$ws.Range("G2").Value = 0.2774390243902439
$ws.Range("H2").Value = 0.6067073170731707
$ws.Range("I2").Value = 0.6128048780487805
$ws.Range("J2").Value = 0.2378048780487805
$ws.Range("K2").Value = 0.4024390243902439
$ws.Range("L2").Value = 0.5335365853658537
$ws.Range("M2").Value = 0.225609756097561
$ws.Range("N2").Value = 0.4481707317073171
$ws.Range("O2").Value = 0.2195121951219512
$ws.Range("P2").Value = 0.4939024390243902
# Row 3
$ws.Range("A3").Value = "BSW"
$ws.Range("B3").Value = 0.4451219512195122
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.4329268292682927
$ws.Range("E3").Value = 0.625
$ws.Range("F3").Value = 0.5701219512195121
$ws.Range("G3").Value = 0.5823170731707317
$ws.Range("H3").Value = 0.5182926829268293
$ws.Range("I3").Value = 0.6280487804878049
$ws.Range("J3").Value = 0.6371951219512195
$ws.Range("K3").Value = 0.5335365853658537
$ws.Range("L3").Value = 0.5945121951219512
$ws.Range("M3").Value = 0.6341463414634146
$ws.Range("N3").Value = 0.6615853658536586
$ws.Range("O3").Value = 0.6097560975609756
$ws.Range("P3").Value = 0.6280487804878049
# Row 4
$ws.Range("A4").Value = "CDU"
$ws.Range("B4").Value = 0.5914634146341463
$ws.Range("C4").Value = 0.4329268292682927
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.4420731707317073
$ws.Range("F4").Value = 0.4939024390243902
$ws.Range("G4").Value = 0.2896341463414634
$ws.Range("H4").Value = 0.6524390243902439
$ws.Range("I4").Value = 0.5792682926829268
$ws.Range("J4").Value = 0.3475609756097561
$ws.Range("K4").Value = 0.5670731707317073
$ws.Range("L4").Value = 0.5579268292682927
$ws.Range("M4").Value = 0.375
$ws.Range("N4").Value = 0.4969512195121951
$ws.Range("O4").Value = 0.399390243902439
$ws.Range("P4").Value = 0.5823170731707317
# Row 5
$ws.Range("A5").Value = "DaGe"
$ws.Range("B5").Value = 0.3170731707317073
$ws.Range("C5").Value = 0.625
$ws.Range("D5").Value = 0.4420731707317073
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.7835365853658537
$ws.Range("G5").Value = 0.7317073170731707
$ws.Range("H5").Value = 0.5396341463414634
$ws.Range("I5").Value = 0.5792682926829268
$ws.Range("J5").Value = 0.7957317073170732
$ws.Range("K5").Value = 0.6676829268292683
$ws.Range("L5").Value = 0.5914634146341463
$ws.Range("M5").Value = 0.7530487804878049
$ws.Range("N5").Value = 0.676829268292683
$ws.Range("O5").Value = 0.8323170731707317
$ws.Range("P5").Value = 0.6554878048780488
# Row 6
$ws.Range("A6").Value = "Die Grünen"
$ws.Range("B6").Value = 0.2347560975609756
$ws.Range("C6").Value = 0.5701219512195121
$ws.Range("D6").Value = 0.4939024390243902
$ws.Range("E6").Value = 0.7835365853658537
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7012195121951219
$ws.Range("H6").Value = 0.5091463414634146
$ws.Range("I6").Value = 0.5335365853658537
$ws.Range("J6").Value = 0.7652439024390244
$ws.Range("K6").Value = 0.7347560975609756
$ws.Range("L6").Value = 0.573170731707317
$ws.Range("M6").Value = 0.75
$ws.Range("N6").Value = 0.6585365853658537
$ws.Range("O6").Value = 0.8414634146341463
$ws.Range("P6").Value = 0.600609756097561
# Row 7
$ws.Range("A7").Value = "Die Partei"
$ws.Range("B7").Value = 0.2774390243902439
$ws.Range("C7").Value = 0.5823170731707317
$ws.Range("D7").Value = 0.2896341463414634
$ws.Range("E7").Value = 0.7317073170731707
$ws.Range("F7").Value = 0.7012195121951219
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.4054878048780488
$ws.Range("I7").Value = 0.5030487804878049
$ws.Range("J7").Value = 0.8536585365853658
$ws.Range("K7").Value = 0.5640243902439024
$ws.Range("L7").Value = 0.4969512195121951
$ws.Range("M7").Value = 0.7896341463414634
$ws.Range("N7").Value = 0.6371951219512195
$ws.Range("O7").Value = 0.7957317073170732
$ws.Range("P7").Value = 0.5945121951219512
# Row 8
$ws.Range("A8").Value = "FDP"
$ws.Range("B8").Value = 0.6067073170731707
$ws.Range("C8").Value = 0.5182926829268293
$ws.Range("D8").Value = 0.6524390243902439
$ws.Range("E8").Value = 0.5396341463414634
$ws.Range("F8").Value = 0.5091463414634146
$ws.Range("G8").Value = 0.4054878048780488
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0.5640243902439024
$ws.Range("J8").Value = 0.4390243902439024
$ws.Range("K8").Value = 0.6128048780487805
$ws.Range("L8").Value = 0.6371951219512195
$ws.Range("M8").Value = 0.4481707317073171
$ws.Range("N8").Value = 0.551829268292683
$ws.Range("O8").Value = 0.4603658536585366
$ws.Range("P8").Value = 0.625
# Row 9
$ws.Range("A9").Value = "Freie Wähler"
$ws.Range("B9").Value = 0.6128048780487805
$ws.Range("C9").Value = 0.6280487804878049
$ws.Range("D9").Value = 0.5792682926829268
$ws.Range("E9").Value = 0.5792682926829268
$ws.Range("F9").Value = 0.5335365853658537
$ws.Range("G9").Value = 0.5030487804878049
$ws.Range("H9").Value = 0.5640243902439024
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0.5091463414634146
$ws.Range("K9").Value = 0.5274390243902439
$ws.Range("L9").Value = 0.6219512195121951
$ws.Range("M9").Value = 0.5182926829268293
$ws.Range("N9").Value = 0.573170731707317
$ws.Range("O9").Value = 0.5182926829268293
$ws.Range("P9").Value = 0.6707317073170732
# Row 10
$ws.Range("A10").Value = "Linke"
$ws.Range("B10").Value = 0.2378048780487805
$ws.Range("C10").Value = 0.6371951219512195
$ws.Range("D10").Value = 0.3475609756097561
$ws.Range("E10").Value = 0.7957317073170732
$ws.Range("F10").Value = 0.7652439024390244
$ws.Range("G10").Value = 0.8536585365853658
$ws.Range("H10").Value = 0.4390243902439024
$ws.Range("I10").Value = 0.5091463414634146
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 0.6128048780487805
$ws.Range("L10").Value = 0.5335365853658537
$ws.Range("M10").Value = 0.8109756097560976
$ws.Range("N10").Value = 0.7012195121951219
$ws.Range("O10").Value = 0.8536585365853658
$ws.Range("P10").Value = 0.6097560975609756
# Row 11
$ws.Range("A11").Value = "SPD"
$ws.Range("B11").Value = 0.4024390243902439
$ws.Range("C11").Value = 0.5335365853658537
$ws.Range("D11").Value = 0.5670731707317073
$ws.Range("E11").Value = 0.6676829268292683
$ws.Range("F11").Value = 0.7347560975609756
$ws.Range("G11").Value = 0.5640243902439024
$ws.Range("H11").Value = 0.6128048780487805
$ws.Range("I11").Value = 0.5274390243902439
$ws.Range("J11").Value = 0.6128048780487805
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.6371951219512195
$ws.Range("M11").Value = 0.6158536585365854
$ws.Range("N11").Value = 0.6615853658536586
$ws.Range("O11").Value = 0.6859756097560976
$ws.Range("P11").Value = 0.6219512195121951
# Row 12
$ws.Range("A12").Value = "SfD"
$ws.Range("B12").Value = 0.5335365853658537
$ws.Range("C12").Value = 0.5945121951219512
$ws.Range("D12").Value = 0.5579268292682927
$ws.Range("E12").Value = 0.5914634146341463
$ws.Range("F12").Value = 0.573170731707317
$ws.Range("G12").Value = 0.4969512195121951
$ws.Range("H12").Value = 0.6371951219512195
$ws.Range("I12").Value = 0.6219512195121951
$ws.Range("J12").Value = 0.5335365853658537
$ws.Range("K12").Value = 0.6371951219512195
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.5487804878048781
$ws.Range("N12").Value = 0.6310975609756098
$ws.Range("O12").Value = 0.5335365853658537
$ws.Range("P12").Value = 0.698170731707317
# Row 13
$ws.Range("A13").Value = "Tierschutzpartei"
$ws.Range("B13").Value = 0.225609756097561
$ws.Range("C13").Value = 0.6341463414634146
$ws.Range("D13").Value = 0.375
$ws.Range("E13").Value = 0.7530487804878049
$ws.Range("F13").Value = 0.75
$ws.Range("G13").Value = 0.7896341463414634
$ws.Range("H13").Value = 0.4481707317073171
$ws.Range("I13").Value = 0.5182926829268293
$ws.Range("J13").Value = 0.8109756097560976
$ws.Range("K13").Value = 0.6158536585365854
$ws.Range("L13").Value = 0.5487804878048781
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 0.6707317073170732
$ws.Range("O13").Value = 0.8079268292682927
$ws.Range("P13").Value = 0.600609756097561
# Row 14
$ws.Range("A14").Value = "UFFBASSE"
$ws.Range("B14").Value = 0.4481707317073171
$ws.Range("C14").Value = 0.6615853658536586
$ws.Range("D14").Value = 0.4969512195121951
$ws.Range("E14").Value = 0.676829268292683
$ws.Range("F14").Value = 0.6585365853658537
$ws.Range("G14").Value = 0.6371951219512195
$ws.Range("H14").Value = 0.551829268292683
$ws.Range("I14").Value = 0.573170731707317
$ws.Range("J14").Value = 0.7012195121951219
$ws.Range("K14").Value = 0.6615853658536586
$ws.Range("L14").Value = 0.6310975609756098
$ws.Range("M14").Value = 0.6707317073170732
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0.6737804878048781
$ws.Range("P14").Value = 0.6646341463414634
# Row 15
$ws.Range("A15").Value = "Volt"
$ws.Range("B15").Value = 0.2195121951219512
$ws.Range("C15").Value = 0.6097560975609756
$ws.Range("D15").Value = 0.399390243902439
$ws.Range("E15").Value = 0.8323170731707317
$ws.Range("F15").Value = 0.8414634146341463
$ws.Range("G15").Value = 0.7957317073170732
$ws.Range("H15").Value = 0.4603658536585366
$ws.Range("I15").Value = 0.5182926829268293
$ws.Range("J15").Value = 0.8536585365853658
$ws.Range("K15").Value = 0.6859756097560976
$ws.Range("L15").Value = 0.5335365853658537
$ws.Range("M15").Value = 0.8079268292682927
$ws.Range("N15").Value = 0.6737804878048781
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = 0.5975609756097561
# Row 16
$ws.Range("A16").Value = "WGD"
$ws.Range("B16").Value = 0.4939024390243902
$ws.Range("C16").Value = 0.6280487804878049
$ws.Range("D16").Value = 0.5823170731707317
$ws.Range("E16").Value = 0.6554878048780488
$ws.Range("F16").Value = 0.600609756097561
$ws.Range("G16").Value = 0.5945121951219512
$ws.Range("H16").Value = 0.625
$ws.Range("I16").Value = 0.6707317073170732
$ws.Range("J16").Value = 0.6097560975609756
$ws.Range("K16").Value = 0.6219512195121951
$ws.Range("L16").Value = 0.698170731707317
$ws.Range("M16").Value = 0.600609756097561
$ws.Range("N16").Value = 0.6646341463414634
$ws.Range("O16").Value = 0.5975609756097561
$ws.Range("P16").Value = 1
